$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns C (NAME), D (SENTENCES) and E (INTELLIGIBILITY_SCORE)
# for rows 2-9. Columns A and B are unchanged.

$ws.Range("C2").Value = "P1_W1_S4"
$ws.Range("D2").Value = "Enjoy the fair weather while in the tropics."
$ws.Range("E2").Value = 0.08163265306122448

$ws.Range("C3").Value = "P1_W1_S1"
$ws.Range("D3").Value = "We picked grapes for wine"
$ws.Range("E3").Value = 0.1333333333333333

$ws.Range("C4").Value = "P1_W1_S2"
$ws.Range("D4").Value = "The ballet is about to begin."
$ws.Range("E4").Value = 0.1764705882352941

$ws.Range("C5").Value = "P1_W1_S3"
$ws.Range("D5").Value = "You're used to being on the field."
$ws.Range("E5").Value = 0.1538461538461539

$ws.Range("C6").Value = "P1_W2_S4"
$ws.Range("D6").Value = "he is capable and willing to make decisions."
$ws.Range("E6").Value = 0.1224489795918367

$ws.Range("C7").Value = "P1_W2_S1"
$ws.Range("D7").Value = "I think I'm getting better."
$ws.Range("E7").Value = 0.125

$ws.Range("C8").Value = "P1_W2_S2"
$ws.Range("D8").Value = "You want him to do well"
$ws.Range("E8").Value = 0.1428571428571428

$ws.Range("C9").Value = "P1_W2_S3"
$ws.Range("D9").Value = "Big muscles are not necessarily strong ones"
$ws.Range("E9").Value = 0.125
